$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.139.10"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("D3").Value = "2.219.73"
$ws.Range("E3").Value = "  -6.07%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.53"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.22"
$ws.Range("E6").Value = "  -9.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.583"
$ws.Range("E7").Value = "  -7.09%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.561"
$ws.Range("E9").Value = "  -8.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.71"
$ws.Range("E10").Value = "  -11.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.34"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0830"
$ws.Range("E12").Value = "  -9.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.77"
$ws.Range("E13").Value = "  -8.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.106"
$ws.Range("E14").Value = "  -3.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.868"
$ws.Range("E15").Value = "  -11.57%  "
$ws.Range("D16").Value = "2.563.44"
$ws.Range("E16").Value = "  -5.84%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.13"
$ws.Range("E17").Value = "  -7.19%  "
$ws.Range("D18").Value = "2.210.19"
$ws.Range("E18").Value = "  -6.38%  "
$ws.Range("D19").Value = "42.974.87"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -9.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -12.36%  "
# Row 23 (special - name swap)
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.32"
$ws.Range("E23").Value = "  -10.70%  "
# Row 24 (special - name swap)
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.20"
$ws.Range("E24").Value = "  -8.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "237.40"
$ws.Range("E25").Value = "  -8.64%  "
$ws.Range("E26").Value = "  -8.18%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("E28").Value = "  -9.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").Value = "  -5.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.28"
$ws.Range("E30").Value = "  -14.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0888"
$ws.Range("E31").Value = "  -7.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.47"
$ws.Range("E32").Value = "  -8.32%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.89"
$ws.Range("E33").Value = "  -9.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "155.04"
$ws.Range("E34").Value = "  -8.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.78"
$ws.Range("E35").Value = "  -5.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.23"
$ws.Range("E36").Value = "  +9.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.02"
$ws.Range("E37").Value = "  +16.56%  "
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.46"
$ws.Range("E39").Value = "  -7.20%  "
$ws.Range("E40").Value = "  -11.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.70"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("E42").Value = "  -8.39%  "
$ws.Range("D43").Value = "1.873.50"
$ws.Range("E43").Value = "  +12.25%  "
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.32"
$ws.Range("E45").Value = "  -4.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.09"
$ws.Range("E46").Value = "  -10.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.207"
$ws.Range("E47").Value = "  -10.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.46"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.07"
$ws.Range("E49").Value = "  -5.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.42"
$ws.Range("E50").Value = "  -13.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.68"
$ws.Range("E51").Value = "  -6.14%  "
